$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '281.07'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '1.30%'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '28.20'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '3.36%'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.042'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '4.01%'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.06484'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '2.30%'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '7.238'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '3.00%'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.376'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '2.09%'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.385'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '1.89%'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.9311'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '4.42%'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1545'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '1.29%'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.06145'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '17.57%'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07547'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '1.42%'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.02899'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '0.21%'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.08988'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '0.45%'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.001608'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '2.42%'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.04439'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '0.97%'
$ws.Range('B17').Value = 'One'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0006371'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '-0.09%'
$ws.Range('B18').Value = 'TigerCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.006047'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '-0.24%'
$ws.Range('B19').Value = 'LEO'
$ws.Range('C19').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.440'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '-0.87%'
$ws.Range('B20').Value = 'BTSEToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '2.232'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '-0.57%'
$ws.Range('B21').Value = 'BitpandaEcosystemToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.3190'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '3.14%'
$ws.Range('B22').Value = 'ProBitToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.1302'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '-2.33%'
$ws.Range('B23').Value = 'MCDex'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.078'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '4.27%'
$ws.Range('B24').Value = 'ZBToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.1547'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '0.55%'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '0.40%'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.004386'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '3.23%'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0001252'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '5.92%'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0001621'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '-1.81%'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '2.65%'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.006649'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '-2.41%'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '-13.49%'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.002024'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '-2.88%'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.01206'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '9.14%'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00005621'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '4.75%'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '25.93%'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.01302'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '-29.74%'
